{"js": "// 1) Fill in the last (previously empty) row of the \"\u6821\u9a8c\u9519\u8bef\" table with the\n//    new error-code entry: 1XX1018 / \u6536\u6b3e\u8d26\u6237\u94f6\u884c\u7c7b\u578b\u9519\u8bef.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// This is the second table in the document body (19 rows: 1 header + 17\n// existing codes + 1 blank trailing row that this edit populates).\nconst errorTable = tables.items[1];\nerrorTable.load(\"rowCount\");\nawait context.sync();\n\nconst lastRowIndex = errorTable.rowCount - 1;\n\nfunction cellOoxml(text) {\n  return `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:hint=\"eastAsia\"/><w:vertAlign w:val=\"baseline\"/><w:lang w:val=\"en-US\" w:eastAsia=\"zh-CN\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/><w:vertAlign w:val=\"baseline\"/><w:lang w:val=\"en-US\" w:eastAsia=\"zh-CN\"/></w:rPr><w:t>${text}</w:t></w:r></w:p></w:body></w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n}\n\nconst codeCell = errorTable.getCell(lastRowIndex, 0);\nconst codeRange = codeCell.body.paragraphs.getFirst().getRange();\ncodeRange.insertOoxml(cellOoxml(\"1XX1018\"), Word.InsertLocation.replace);\nawait context.sync();\n\nconst descCell = errorTable.getCell(lastRowIndex, 1);\nconst descRange = descCell.body.paragraphs.getFirst().getRange();\ndescRange.insertOoxml(cellOoxml(\"\u6536\u6b3e\u8d26\u6237\u94f6\u884c\u7c7b\u578b\u9519\u8bef\"), Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Drop the stray \"_GoBack\" bookmark left over from the last cursor\n//    position (Word auto-manages this bookmark; it is not meaningful\n//    content and is removed here).\nif (context.document.bookmarks.exists(\"_GoBack\")) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 3) Recolor the theme's \"Light 1\" swatch (window/background color) from\n//    white to a light green (CCE8CF).\n//    NOTE: the Word JavaScript API does not expose the document's DrawingML\n//    color theme (word/theme/theme1.xml) for editing - there is no\n//    Word.Document member for the theme color scheme, and\n//    Range/Body.insertOoxml only ever materializes the /word/document.xml\n//    part of a supplied OOXML package, so a theme part cannot be pushed\n//    through it either. This part of the change has no Office.js-reachable\n//    equivalent, so it is intentionally left for the COM/VBA object model\n//    (Document.DocumentTheme.ThemeColorScheme), which does expose it.\n", "ps1": "# Word COM interop script - applies the same three changes as edit.js:\n#   1) Fill in the last (previously empty) row of the \"\u6821\u9a8c\u9519\u8bef\" error-code\n#      table with 1XX1018 / \u6536\u6b3e\u8d26\u6237\u94f6\u884c\u7c7b\u578b\u9519\u8bef.\n#   2) Remove the stray \"_GoBack\" bookmark left at the cursor's last edit\n#      position.\n#   3) Recolor the document theme's \"Light 1\" swatch from white to a light\n#      green (CCE8CF).\n\n$d = $word.ActiveDocument\n\n# --- 1) Populate the trailing blank row of the error-code table ----------\n# The table is the 2nd table in the document (1-based Tables collection)\n# and its last row is still empty (code/description both blank).\n$table = $d.Tables.Item(2)\n$lastRow = $table.Rows.Count\n\n$codeCell = $table.Cell($lastRow, 1)\n$codeCell.Range.Text = \"1XX1018\"\n\n$descCell = $table.Cell($lastRow, 2)\n$descCell.Range.Text = \"\u6536\u6b3e\u8d26\u6237\u94f6\u884c\u7c7b\u578b\u9519\u8bef\"\n\n# --- 2) Remove the leftover _GoBack bookmark ------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- 3) Recolor theme \"Light 1\" (window/background) from white to CCE8CF -\n$themeColors = $d.DocumentTheme.ThemeColorScheme\n# Scheme order: 1=dk1 2=lt1 3=dk2 4=lt2 5-10=accent1-6 11=hlink 12=folHlink\n$lt1 = $themeColors.Colors(2)\n$lt1.RGB = 13625548   # 0xCCE8CF packed as BGR (R | G<<8 | B<<16) = RGB(204,232,207)\n"}
